$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 282 (shifts existing rows 282-296 down to 283-297)
$ws.Rows("282:282").Insert()

# Populate the newly inserted row with the new Jengibre price record
$ws.Range("A282").Value = 10
$ws.Range("B282").Value = "Vega Modelo de Temuco"
$ws.Range("C282").Value = "La Araucanía"
$ws.Range("D282").Value = 45075
$ws.Range("E282").Value = 9
$ws.Range("F282").Value = 100114007
$ws.Range("G282").Value = "Jengibre"
$ws.Range("H282").Value = "Sin especificar"
$ws.Range("I282").Value = "Primera"
$ws.Range("J282").Value = 15
$ws.Range("K282").Value = 24000
$ws.Range("L282").Value = 24000
$ws.Range("M282").Value = 24000
$ws.Range("N282").Value = "`$/caja 13 kilos"
$ws.Range("O282").Value = "Perú"
$ws.Range("P282").Value = 1846
$ws.Range("Q282").Value = 13
$ws.Range("R282").Value = "Hortaliza"
